# Auto-generated edit script
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Black River Falls Wisconsin")
$ws.Range("L7").ClearContents()
$ws.Range("E8").Value = 0.0776
$ws.Range("E9").Value = 0.0776
$ws.Range("E10").Value = 0.0776
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 0
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = 0
$ws.Range("Q10").Value = 0
$ws.Range("R10").Value = 0
$ws.Range("S10").Value = 0
$ws.Range("T10").Value = 0
$ws.Range("U10").Value = 0
$ws.Range("V10").Value = 0
$ws.Range("W10").Value = 0

$ws = $wb.Worksheets.Item("Cassville Missouri")
$ws.Range("E5").Value = 0.0776
$ws.Range("E6").Value = 0.0776
$ws.Range("E7").Value = 0.0776
$ws.Range("L7").Value = 0.0159
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 0
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 0
$ws.Range("R7").Value = 0
$ws.Range("S7").Value = 0
$ws.Range("T7").Value = 0
$ws.Range("U7").Value = 0
$ws.Range("V7").Value = 0
$ws.Range("W7").Value = 0

$ws = $wb.Worksheets.Item("Eaglepass Texas")
$ws.Range("L4").Value = 0
$ws.Range("E7").Value = 0.0776
$ws.Range("E8").Value = 0.0776
$ws.Range("E9").Value = 0.0776

$ws = $wb.Worksheets.Item("Faridabad India")
$ws.Range("L7").ClearContents()
$ws.Range("E8").Value = 0.0776
$ws.Range("E9").Value = 0.0776
$ws.Range("E10").Value = 0.0776

$ws = $wb.Worksheets.Item("Fort Wayne Indiana")
$ws.Range("L7").ClearContents()
$ws.Range("E8").Value = 0.0776
$ws.Range("E9").Value = 0.0776

$ws = $wb.Worksheets.Item("Juarez Casa I")
$ws.Range("L4").Value = 0.0204
$ws.Range("L7").ClearContents()
$ws.Range("E8").Value = 0.0776
$ws.Range("E9").Value = 0.0776
$ws.Range("E10").Value = 0.0776
$ws.Range("J10").Value = 0.0379
$ws.Range("L10").Value = 0.0082

$ws = $wb.Worksheets.Item("Juarez Casa II")
$ws.Range("E7").Value = 0.0776
$ws.Range("E8").Value = 0.0776
$ws.Range("E9").Value = 0.0776
$ws.Range("I9").Value = 0.007
$ws.Range("J9").Value = 0.0212
$ws.Range("K9").Value = 0.0199
$ws.Range("L9").Value = 0.0094

$ws = $wb.Worksheets.Item("Juarez Casa SS")
$ws.Range("L5").ClearContents()
$ws.Range("E6").Value = 0.0776
$ws.Range("E7").Value = 0.0776
$ws.Range("E8").Value = 0.0776

$ws = $wb.Worksheets.Item("Juarez FCDM")
$ws.Range("E8").Value = 0.0776
$ws.Range("E9").Value = 0.0776
$ws.Range("E10").Value = 0.0776
$ws.Range("G10").Value = 0.0154
$ws.Range("J10").Value = 0.0291
$ws.Range("K10").Value = 0.0104
$ws.Range("L10").Value = 0.0157

$ws = $wb.Worksheets.Item("Juarez MEJ II")
$ws.Range("L4").Value = 0.0333
$ws.Range("E7").Value = 0.0776
$ws.Range("E8").Value = 0.0776
$ws.Range("E9").Value = 0.0776
$ws.Range("I9").Value = 0.0229
$ws.Range("J9").Value = 0.0771
$ws.Range("L9").Value = 0.0264

$ws = $wb.Worksheets.Item("Mcallen Texas")
$ws.Range("L4").Value = 0.0588
$ws.Range("L7").ClearContents()
$ws.Range("E8").Value = 0.0776
$ws.Range("E9").Value = 0.0776
$ws.Range("E10").Value = 0.0776
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 0
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = 0
$ws.Range("Q10").Value = 0
$ws.Range("R10").Value = 0
$ws.Range("S10").Value = 0
$ws.Range("T10").Value = 0
$ws.Range("U10").Value = 0
$ws.Range("V10").Value = 0
$ws.Range("W10").Value = 0

$ws = $wb.Worksheets.Item("Piedras Negras Fasco Mexico")
$ws.Range("H4").Value = 0.0169
$ws.Range("J4").Value = 0.0166
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0
$ws.Range("L7").ClearContents()
$ws.Range("E8").Value = 0.0776
$ws.Range("E9").Value = 0.0776
$ws.Range("E10").Value = 0.0776
$ws.Range("H10").Value = 0.0197
$ws.Range("I10").Value = 0.0221
$ws.Range("J10").Value = 0.0619
$ws.Range("K10").Value = 0.0277
$ws.Range("L10").Value = 0.0214

$ws = $wb.Worksheets.Item("Piedras Negras Jakel Mexico")
$ws.Range("G4").Value = 0.1429
$ws.Range("J4").Value = 0.1429
$ws.Range("L7").ClearContents()
$ws.Range("E8").Value = 0.0776
$ws.Range("E9").Value = 0.0776
$ws.Range("E10").Value = 0.0776
$ws.Range("G10").Value = 0.1096
$ws.Range("H10").Value = 0.1935
$ws.Range("I10").Value = 0.1081
$ws.Range("J10").Value = 0.4094
$ws.Range("K10").Value = 0.0909
$ws.Range("L10").Value = 0.0638

$ws = $wb.Worksheets.Item("Reynosa Mexico")
$ws.Range("E5").Value = 0.0776
$ws.Range("E6").Value = 0.0776

$ws = $wb.Worksheets.Item("Tipp City Ohio")
$ws.Range("L4").Value = 0.0513

$ws = $wb.Worksheets.Item("Milwaukee Pmc Hq Wisconsin")
$ws.Range("L3").ClearContents()

$ws = $wb.Worksheets.Item("Reynosa II")
$ws.Range("E2").Value = 0.0776
$ws.Range("E3").Value = 0.0776

$ws = $wb.Worksheets.Item("Sao Paulo Brazil")
$ws.Range("E2").Value = 0.0776
$ws.Range("E3").Value = 0.0776
